$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/3/2024  Through  6/9/2024"

# --- Crime Complaints table updates (rows 15-31) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = -25
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 50

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -66.666666666666
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -38.888888888888
$ws.Range("M16").Value = -27.868852459016
$ws.Range("N16").Value = -86.666666666666

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -62.5
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -35
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 89
$ws.Range("K17").Value = 17.977528089887
$ws.Range("L17").Value = 10.526315789473
$ws.Range("M17").Value = 59.090909090909
$ws.Range("N17").Value = 9.375

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 63
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = -4.545454545454
$ws.Range("L18").Value = -12.5
$ws.Range("M18").Value = 125
$ws.Range("N18").Value = -58

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -38.181818181818
$ws.Range("I19").Value = 171
$ws.Range("J19").Value = 243
$ws.Range("K19").Value = -29.629629629629
$ws.Range("L19").Value = -50
$ws.Range("M19").Value = 67.647058823529
$ws.Range("N19").Value = -11.855670103092

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -37.5
$ws.Range("I20").Value = 25
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = 4.166666666666
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -84.076433121019

# Row 21
$ws.Range("C21").Value = 14
$ws.Range("E21").Value = -53.333333333333
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -36.842105263157
$ws.Range("I21").Value = 414
$ws.Range("J21").Value = 483
$ws.Range("K21").Value = -14.285714285714
$ws.Range("L21").Value = -32.792207792207
$ws.Range("M21").Value = 45.263157894736
$ws.Range("N21").Value = -55.769230769230

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("H22").Value = 100
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -11.111111111111
$ws.Range("L22").Value = -52.941176470588
$ws.Range("M22").Value = 100

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 33.333333333333
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = 18.75
$ws.Range("I23").Value = 97
$ws.Range("J23").Value = 63
$ws.Range("K23").Value = 53.968253968254
$ws.Range("L23").Value = 19.753086419753
$ws.Range("M23").Value = 38.571428571428

# Row 24
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -15.151515151515
$ws.Range("F24").Value = 129
$ws.Range("G24").Value = 156
$ws.Range("H24").Value = -17.307692307692
$ws.Range("I24").Value = 563
$ws.Range("J24").Value = 527
$ws.Range("K24").Value = 6.831119544592
$ws.Range("L24").Value = -48.062730627306
$ws.Range("M24").Value = 80.448717948717

# Row 25
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = 10.526315789473
$ws.Range("F25").Value = 95
$ws.Range("G25").Value = 85
$ws.Range("H25").Value = 11.764705882352
$ws.Range("I25").Value = 371
$ws.Range("J25").Value = 281
$ws.Range("K25").Value = 32.028469750889
$ws.Range("L25").Value = -58.173618940248

# Row 26
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 52
$ws.Range("H26").Value = -42.307692307692
$ws.Range("I26").Value = 194
$ws.Range("J26").Value = 191
$ws.Range("K26").Value = 1.570680628272
$ws.Range("L26").Value = -4.901960784313
$ws.Range("M26").Value = 37.588652482269

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 12
$ws.Range("K27").Value = 9.090909090909
$ws.Range("L27").Value = -14.285714285714

# Row 28
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -40
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 26
$ws.Range("K28").Value = -38.461538461538
$ws.Range("L28").Value = -5.882352941176

# Row 31
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = 200
$ws.Range("L31").Value = 12.5
